$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the old title row and clear old hyperlinks ---
$ws.Range("A1:G1").UnMerge()
$ws.Hyperlinks.Delete()

# --- Delete the now-unused trailing columns F:G ---
$ws.Columns("F:G").Delete()

# --- Write the new header + data rows ---
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "surname"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "interest"

$ws.Range("A2").Value = "John"
$ws.Range("B2").Value = "Smith"
$ws.Range("C2").Value = "ypddjuio@yomail.info"
$ws.Range("D2").Value = "meditation"

$ws.Range("A3").Value = "Marry"
$ws.Range("B3").Value = "Smith"
$ws.Range("C3").Value = "kawnlyiw@supere.ml"
$ws.Range("D3").Value = "nasa"

$ws.Range("A4").Value = "Sim"
$ws.Range("B4").Value = "Kann"
$ws.Range("C4").Value = "pythonprocourse2@gmail.com"
$ws.Range("D4").Value = "nasa"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ypddjuio@yomail.info", "", "", "ypddjuio@yomail.info")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:kawnlyiw@supere.ml", "", "", "kawnlyiw@supere.ml")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:pythonprocourse2@gmail.com", "", "", "pythonprocourse2@gmail.com")

Write-Host "done"
